$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All Price/Volume cells are stored as text (e.g. thousand-dot-separated
# prices, percentages padded with spaces) - force text format so Excel
# does not auto-convert the assigned strings into numbers.

# --- Update Price (D) and Volume(1h) (E) columns for changed rows ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.862.92"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.15%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.905.48"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.05%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.54"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.66"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.48%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.902.43"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.54%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.99"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.79%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.147"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.60%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.17%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.90%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.27"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.41%  "

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.16%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.387.23"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.847.20"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.19%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.908.35"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.38%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.56%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "428.77"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.94%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.41%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.650"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.12%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.61%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.66"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.98"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.16"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -7.03%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.02"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.25%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +11.68%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.96"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.50%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.01"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.92%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.20%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.17%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.52"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.33%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.956"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.73%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.19%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.79"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.80%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -5.77%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.89"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.34%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.14"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.42%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.266"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.701.82"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.73%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.13%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "132.57"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.91%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "347.48"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.54%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.76%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.56"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.75%  "

# --- Rows 41/42: Kaspa and Arweave swap places (rank in column A stays the same) ---
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.114"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.65%  "

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Arweave"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.98"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.61%  "
